$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "332.48"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.30%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.64"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.34%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.736"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.02%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08111"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.18%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.684"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.02%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.488"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.71%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.982"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.78%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.000"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.72%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9288"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.57%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1277"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.49%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1959"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.18%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.778"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "15.65%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09167"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.02%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03760"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "8.60%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.1050"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "9.21%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001294"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.27%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006304"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.75%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.369"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.09%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3496"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.93%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1367"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.15%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "7.13%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04434"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.21%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.01%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004405"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.83%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.93%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02823"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "11.92%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05604"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "7.19%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007514"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.01%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009816"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "9.04%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1424"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.38%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.97%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01186"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "18.81%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006780"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.28%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.28%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003064"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "6.67%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002272"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "26.13%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002094"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.28%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001994"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.28%"
